# Orchestrator credentials are now stored in Windows Credential Manager (WCM).
#
# This reshapes the "Settings" sheet so the orchestrator-related rows use the
# new naming/values, removes the now-obsolete username/password rows, and
# switches the active sheet/tab from "Assets" back to "Settings".

$wb = $excel.ActiveWorkbook
$settings = $wb.Worksheets.Item("Settings")

# Rows 2-5: replace the old Orch_* rows with the new Orchestrator* rows,
# each now carrying Name / Value / Description.
$settings.Range("A2").Value = "OrchestratorURL"
$settings.Range("B2").Value = "https://demo.uipath.com"
$settings.Range("C2").Value = "The URL of your orchestrator server. This property is used only if you are using a Queue to store your Transaction Items."

$settings.Range("A3").Value = "OrchestratorCredentialName"
$settings.Range("B3").Value = "demo.uipath.com_credentials"
$settings.Range("C3").Value = "The name of Orchestrator credentials. This should be stored in Windows Credential manager. This property is used only if you are using a Queue to store your Transaction Items."

$settings.Range("A4").Value = "OrchestratorTenancyName"
$settings.Range("B4").Value = "fantastic"
$settings.Range("C4").Value = "The name of the Orchestrator tenant.  This property is used only if you are using a Queue to store your Transaction Items."

$settings.Range("A5").Value = "OrchestratorQueueName"
$settings.Range("B5").Value = "KibanaDemoQueue"
$settings.Range("C5").Value = "Orchestrator Queue Name. Be sure to match this name with the one from the server."

# The old sheet had the credential/username/password rows trailing at 6:8 -
# those no longer exist now that credentials live in WCM, so drop the rows
# (shrinking the sheet dimension down to A1:C5).
$settings.Rows("6:8").Delete()

# Make "Settings" the active tab again (it had drifted to "Assets"), with
# the selection resting on A2.
$settings.Activate()
$settings.Range("A2").Select()
